$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($r = 1; $r -le 14; $r++) {
    if ($r -eq 14) {
        $val = 0.2
    } else {
        $off = $r - 5
        $val = $off * $off * 0.001
    }
    for ($c = 1; $c -le 8; $c++) {
        $ws.Cells.Item($r, $c).Value = $val
    }
}

$ws.Range("J1:K14").Select()
